# Added optional "enabled" argument column to the Properties sheet,
# and updated the active cell selections on both sheets.

$wb = $excel.ActiveWorkbook
$wsParameters = $wb.Worksheets.Item("Parameters")
$wsProperties = $wb.Worksheets.Item("Properties")

# --- Properties sheet: add new "enabled" column (D) ---
$wsProperties.Range("D1").Value = "enabled"

$lastRow = 21
for ($r = 2; $r -le $lastRow; $r++) {
    $wsProperties.Cells.Item($r, 4).Value = $true
}

# --- Update selections to match the new state ---
# Parameters sheet ends up with D2 selected (but is not the active tab).
$wsParameters.Range("D2").Select()

# Properties sheet is (re)activated and ends up with F11 selected.
$wsProperties.Activate()
$wsProperties.Range("F11").Select()
